# Insert a new price-report row at row 68 (pushing the existing rows 68-123
# down to 69-124, growing the used range from A1:R123 to A1:R124), then
# populate the newly inserted row with the new weekly record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 68..123 down by one row.
$ws.Rows.Item(68).Insert()

# Fill in the new row 68 with the new data point.
$ws.Cells.Item(68, 1).Value = 11
$ws.Cells.Item(68, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(68, 3).Value = "Bíobío"
$ws.Cells.Item(68, 4).Value = 44705
$ws.Cells.Item(68, 5).Value = 8
$ws.Cells.Item(68, 6).Value = 100112032
$ws.Cells.Item(68, 7).Value = "Zapallo italiano"
$ws.Cells.Item(68, 8).Value = "Sin especificar"
$ws.Cells.Item(68, 9).Value = "Primera"
$ws.Cells.Item(68, 10).Value = 150
$ws.Cells.Item(68, 11).Value = 18000
$ws.Cells.Item(68, 12).Value = 19000
$ws.Cells.Item(68, 13).Value = 18333
$ws.Cells.Item(68, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(68, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(68, 16).Value = 367
$ws.Cells.Item(68, 17).Value = 50
$ws.Cells.Item(68, 18).Value = "Hortaliza"
